$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("G2").Value = 1.57
$ws.Range("H2").Value = 3.6
$ws.Range("I2").Value = 6.25
$ws.Range("J2").Value = 1.08
$ws.Range("K2").Value = 8
$ws.Range("N2").Value = 2.25
$ws.Range("O2").Value = 1.62
$ws.Range("R2").Value = 2.2
$ws.Range("S2").Value = 1.62
$ws.Range("U2").Value = 6.5
$ws.Range("V2").Value = 9
$ws.Range("W2").Value = 11
$ws.Range("Z2").Value = 8
$ws.Range("AA2").Value = 7
$ws.Range("AC2").Value = 81
$ws.Range("AE2").Value = 13
$ws.Range("AF2").Value = 29
$ws.Range("AG2").Value = 21
$ws.Range("AI2").Value = 51

# Row 3 - populate previously empty odds cells
$ws.Range("G3").Value = 2.3
$ws.Range("H3").Value = 3.55
$ws.Range("I3").Value = 2.82
$ws.Range("J3").Value = 1.04
$ws.Range("K3").Value = 9
$ws.Range("L3").Value = 1.22
$ws.Range("M3").Value = 4
$ws.Range("N3").Value = 1.65
$ws.Range("O3").Value = 2.15
$ws.Range("P3").Value = 1.35
$ws.Range("Q3").Value = 3.05
$ws.Range("R3").Value = 1.55
$ws.Range("S3").Value = 2.3
$ws.Range("T3").Value = 9.75
$ws.Range("U3").Value = 14
$ws.Range("V3").Value = 9.5
$ws.Range("W3").Value = 26
$ws.Range("X3").Value = 18
$ws.Range("Y3").Value = 24
$ws.Range("Z3").Value = 9
$ws.Range("AA3").Value = 7.3
$ws.Range("AB3").Value = 12.5
$ws.Range("AC3").Value = 50
$ws.Range("AD3").Value = 300
$ws.Range("AE3").Value = 10.75
$ws.Range("AF3").Value = 17.5
$ws.Range("AG3").Value = 10.75
$ws.Range("AH3").Value = 37
$ws.Range("AI3").Value = 23
$ws.Range("AJ3").Value = 28

# Row 4 updates
$ws.Range("G4").Value = 2.35
$ws.Range("H4").Value = 3
$ws.Range("I4").Value = 3.25
$ws.Range("K4").Value = 6.5
$ws.Range("N4").Value = 2.5
$ws.Range("O4").Value = 1.5
$ws.Range("P4").Value = 1.53
$ws.Range("Q4").Value = 2.38
$ws.Range("R4").Value = 2.1
$ws.Range("S4").Value = 1.67
$ws.Range("T4").Value = 6.5
$ws.Range("U4").Value = 10
$ws.Range("V4").Value = 10
$ws.Range("W4").Value = 23
$ws.Range("X4").Value = 23
$ws.Range("Y4").Value = 41
$ws.Range("Z4").Value = 6.5
$ws.Range("AA4").Value = 6
$ws.Range("AB4").Value = 19
$ws.Range("AC4").Value = 67
$ws.Range("AE4").Value = 7.5
$ws.Range("AF4").Value = 15
$ws.Range("AG4").Value = 12
$ws.Range("AH4").Value = 34
$ws.Range("AI4").Value = 29
$ws.Range("AJ4").Value = 41

# Row 5 updates
$ws.Range("G5").Value = 2.1
$ws.Range("H5").Value = 3.1
$ws.Range("I5").Value = 3.7
$ws.Range("K5").Value = 7
$ws.Range("M5").Value = 2.63
$ws.Range("N5").Value = 2.35
$ws.Range("P5").Value = 1.5
$ws.Range("Q5").Value = 2.5
$ws.Range("R5").Value = 2
$ws.Range("S5").Value = 1.73
$ws.Range("T5").Value = 6
$ws.Range("U5").Value = 9
$ws.Range("V5").Value = 9.5
$ws.Range("W5").Value = 19
$ws.Range("Y5").Value = 34
$ws.Range("Z5").Value = 7
$ws.Range("AA5").Value = 6
$ws.Range("AB5").Value = 17
$ws.Range("AC5").Value = 67
$ws.Range("AD5").Value = 1250
$ws.Range("AE5").Value = 9
$ws.Range("AF5").Value = 17
$ws.Range("AG5").Value = 13
$ws.Range("AH5").Value = 41
$ws.Range("AI5").Value = 34
$ws.Range("AJ5").Value = 41

# Row 9 - time change
$ws.Range("C9").Value = "14:30"

# Row 12 - populate previously empty odds cells (J12, K12, AD12 stay empty)
$ws.Range("G12").Value = 2.27
$ws.Range("H12").Value = 2.95
$ws.Range("I12").Value = 3.2
$ws.Range("L12").Value = 1.47
$ws.Range("M12").Value = 2.32
$ws.Range("N12").Value = 2.37
$ws.Range("O12").Value = 1.45
$ws.Range("P12").Value = 1.52
$ws.Range("Q12").Value = 2.2
$ws.Range("R12").Value = 2.02
$ws.Range("S12").Value = 1.62
$ws.Range("T12").Value = 6
$ws.Range("U12").Value = 9.75
$ws.Range("V12").Value = 9.75
$ws.Range("W12").Value = 23
$ws.Range("X12").Value = 22
$ws.Range("Y12").Value = 45
$ws.Range("Z12").Value = 6.5
$ws.Range("AA12").Value = 5.9
$ws.Range("AB12").Value = 18.5
$ws.Range("AC12").Value = 120
$ws.Range("AE12").Value = 7.2
$ws.Range("AF12").Value = 14.5
$ws.Range("AG12").Value = 12
$ws.Range("AH12").Value = 45
$ws.Range("AI12").Value = 37
$ws.Range("AJ12").Value = 55

# Row 13 - populate previously empty odds cells (J13, K13, AD13 stay empty)
$ws.Range("G13").Value = 2.95
$ws.Range("H13").Value = 2.85
$ws.Range("I13").Value = 2.5
$ws.Range("L13").Value = 1.47
$ws.Range("M13").Value = 2.35
$ws.Range("N13").Value = 2.32
$ws.Range("O13").Value = 1.47
$ws.Range("P13").Value = 1.5
$ws.Range("Q13").Value = 2.25
$ws.Range("R13").Value = 1.93
$ws.Range("S13").Value = 1.7
$ws.Range("T13").Value = 7
$ws.Range("U13").Value = 13.5
$ws.Range("V13").Value = 11.25
$ws.Range("W13").Value = 40
$ws.Range("X13").Value = 32
$ws.Range("Y13").Value = 50
$ws.Range("Z13").Value = 6.5
$ws.Range("AA13").Value = 5.7
$ws.Range("AB13").Value = 16.5
$ws.Range("AC13").Value = 100
$ws.Range("AE13").Value = 6.6
$ws.Range("AF13").Value = 11.25
$ws.Range("AG13").Value = 9.75
$ws.Range("AH13").Value = 28
$ws.Range("AI13").Value = 24
$ws.Range("AJ13").Value = 40

# Row 16 - populate previously empty odds cells (all columns)
$ws.Range("G16").Value = 2.35
$ws.Range("H16").Value = 3.6
$ws.Range("I16").Value = 2.67
$ws.Range("J16").Value = 1.03
$ws.Range("K16").Value = 9.25
$ws.Range("L16").Value = 1.16
$ws.Range("M16").Value = 4.6
$ws.Range("N16").Value = 1.5
$ws.Range("O16").Value = 2.42
$ws.Range("P16").Value = 1.29
$ws.Range("Q16").Value = 3.3
$ws.Range("R16").Value = 1.44
$ws.Range("S16").Value = 2.6
$ws.Range("T16").Value = 12.5
$ws.Range("U16").Value = 15
$ws.Range("V16").Value = 9.25
$ws.Range("W16").Value = 26
$ws.Range("X16").Value = 16.5
$ws.Range("Y16").Value = 19
$ws.Range("Z16").Value = 9.25
$ws.Range("AA16").Value = 7.6
$ws.Range("AB16").Value = 10.75
$ws.Range("AC16").Value = 32
$ws.Range("AD16").Value = 175
$ws.Range("AE16").Value = 13.5
$ws.Range("AF16").Value = 17.5
$ws.Range("AG16").Value = 10
$ws.Range("AH16").Value = 32
$ws.Range("AI16").Value = 19
$ws.Range("AJ16").Value = 21
